{"js": "// Journal entry proof-read pass: insert a handful of missing words and\n// carry the \"_GoBack\" (last-edit-position) bookmark forward to the point\n// of the final insertion, the way Word leaves it after an editing session.\nconst body = context.document.body;\n\n// 1) \"Empire of Sun\" -> \"Empire of the Sun\"\nconst hit1 = body.search(\"Sun which was first released\", { matchCase: true });\nhit1.load(\"items\");\nawait context.sync();\nhit1.items[0].insertText(\"the \", \"Before\");\nawait context.sync();\n\n// 2) \"born with silver spoon\" -> \"born with a silver spoon\"\nconst hit2 = body.search(\"silver spoon and had been living\", { matchCase: true });\nhit2.load(\"items\");\nawait context.sync();\nhit2.items[0].insertText(\"a \", \"Before\");\nawait context.sync();\n\n// 3) \"living in China with his family\" -> \"living in China since he was born with his family\"\nconst hit3 = body.search(\"with his family who were colonists\", { matchCase: true });\nhit3.load(\"items\");\nawait context.sync();\nhit3.items[0].insertText(\"since he was born \", \"Before\");\nawait context.sync();\n\n// 4) \"he knew some Japanese war planes\" -> \"he knew about some Japanese war planes\"\nconst hit4 = body.search(\"some Japanese war planes\", { matchCase: true });\nhit4.load(\"items\");\nawait context.sync();\nhit4.items[0].insertText(\"about \", \"Before\");\nawait context.sync();\n\n// 5) Move the \"_GoBack\" bookmark from the end of the document to right before\n//    \"some Japanese war planes\" in the movie paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst hit5 = body.search(\"some Japanese war planes\", { matchCase: true });\nhit5.load(\"items\");\nawait context.sync();\nconst bookmarkAnchor = hit5.items[0].getRange(\"Start\");\nbookmarkAnchor.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Journal entry proof-read pass: insert a handful of missing words and\n# carry the \"_GoBack\" (last-edit-position) bookmark forward to the point\n# of the final insertion, the way Word leaves it after an editing session.\n\n$d = $word.ActiveDocument\n$wdCollapseStart = 1\n\n# 1) \"Empire of Sun\" -> \"Empire of the Sun\"\n$r1 = $d.Content\n$r1.Find.Execute(\"Sun which was first released\")\n$r1.Collapse($wdCollapseStart)\n$r1.InsertBefore(\"the \")\n\n# 2) \"born with silver spoon\" -> \"born with a silver spoon\"\n$r2 = $d.Content\n$r2.Find.Execute(\"silver spoon and had been living\")\n$r2.Collapse($wdCollapseStart)\n$r2.InsertBefore(\"a \")\n\n# 3) \"living in China with his family\" -> \"living in China since he was born with his family\"\n$r3 = $d.Content\n$r3.Find.Execute(\"with his family who were colonists\")\n$r3.Collapse($wdCollapseStart)\n$r3.InsertBefore(\"since he was born \")\n\n# 4) \"he knew some Japanese war planes\" -> \"he knew about some Japanese war planes\"\n$r4 = $d.Content\n$r4.Find.Execute(\"some Japanese war planes\")\n$r4.Collapse($wdCollapseStart)\n$r4.InsertBefore(\"about \")\n\n# 5) Move the \"_GoBack\" bookmark from the end of the document to right\n#    before \"some Japanese war planes\" (the last edit made above).\n$bms = $d.Bookmarks\nif ($bms.Exists(\"_GoBack\")) {\n  $bms.Item(\"_GoBack\").Delete()\n}\n$r5 = $d.Content\n$r5.Find.Execute(\"some Japanese war planes\")\n$r5.Collapse($wdCollapseStart)\n$bms.Add(\"_GoBack\", $r5)\n"}
